$wb = $excel.ActiveWorkbook

# --- "Repayment schedule" sheet: insert a new (blank) column before column N ---
# This is the RBI "Variable Instalments" column that was added to the
# repayment schedule. Inserting pushes the former N/O/P columns
# ("Late" / "heading" / "Outstanding") one place to the right, becoming
# O/P/Q, and leaves a blank column N that inherits column M's formatting.
$schedule = $wb.Worksheets.Item("Repayment schedule")
$schedule.Columns("N").Insert()

# Make "Repayment schedule" the active sheet/tab (it takes over from
# "Transactions", which was previously active) and set its new selection.
$schedule.Activate()
$schedule.Range("R7").Select()
